$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 4).Value = '64.055.73'
$ws.Cells.Item(2, 5).Value = '  +0.51%  '
$ws.Cells.Item(3, 4).Value = '2.777.90'
$ws.Cells.Item(3, 5).Value = '  +1.55%  '
$st = $ws.Cells.Item(4, 4).Style
$ws.Cells.Item(4, 4).NumberFormat = '@'
$ws.Cells.Item(4, 4).Value = '1.00'
$ws.Cells.Item(4, 4).Style = $st
$ws.Cells.Item(4, 5).Value = '  +0.09%  '
$st = $ws.Cells.Item(5, 4).Style
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '587.25'
$ws.Cells.Item(5, 4).Style = $st
$ws.Cells.Item(5, 5).Value = '  -0.48%  '
$st = $ws.Cells.Item(6, 4).Style
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '160.18'
$ws.Cells.Item(6, 4).Style = $st
$ws.Cells.Item(6, 5).Value = '  +5.62%  '
$st = $ws.Cells.Item(7, 4).Style
$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '0.998'
$ws.Cells.Item(7, 4).Style = $st
$ws.Cells.Item(7, 5).Value = '  +0.36%  '
$st = $ws.Cells.Item(8, 4).Style
$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '0.617'
$ws.Cells.Item(8, 4).Style = $st
$ws.Cells.Item(8, 5).Value = '  +0.86%  '
$st = $ws.Cells.Item(9, 4).Style
$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '0.113'
$ws.Cells.Item(9, 4).Style = $st
$ws.Cells.Item(9, 5).Value = '  -0.25%  '
$st = $ws.Cells.Item(10, 4).Style
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '6.03'
$ws.Cells.Item(10, 4).Style = $st
$ws.Cells.Item(10, 5).Value = '  -11.04%  '
$st = $ws.Cells.Item(11, 4).Style
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '0.396'
$ws.Cells.Item(11, 4).Style = $st
$ws.Cells.Item(11, 5).Value = '  +1.69%  '
$st = $ws.Cells.Item(12, 4).Style
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '0.159'
$ws.Cells.Item(12, 4).Style = $st
$ws.Cells.Item(12, 5).Value = '  +0.60%  '
$ws.Cells.Item(13, 4).Value = '3.265.08'
$ws.Cells.Item(13, 5).Value = '  +1.74%  '
$st = $ws.Cells.Item(14, 4).Style
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '27.24'
$ws.Cells.Item(14, 4).Style = $st
$ws.Cells.Item(14, 5).Value = '  +1.47%  '
$ws.Cells.Item(15, 4).Value = '63.947.21'
$ws.Cells.Item(15, 5).Value = '  +0.61%  '
$st = $ws.Cells.Item(16, 4).Style
$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '0.0000158'
$ws.Cells.Item(16, 4).Style = $st
$ws.Cells.Item(16, 5).Value = '  +3.62%  '
$ws.Cells.Item(17, 4).Value = '2.780.94'
$ws.Cells.Item(17, 5).Value = '  +1.31%  '
$st = $ws.Cells.Item(18, 4).Style
$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '12.37'
$ws.Cells.Item(18, 4).Style = $st
$ws.Cells.Item(18, 5).Value = '  +2.64%  '
$st = $ws.Cells.Item(19, 4).Style
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '5.00'
$ws.Cells.Item(19, 4).Style = $st
$ws.Cells.Item(19, 5).Value = '  +2.46%  '
$st = $ws.Cells.Item(20, 4).Style
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '364.44'
$ws.Cells.Item(20, 4).Style = $st
$ws.Cells.Item(20, 5).Value = '  -0.31%  '
$st = $ws.Cells.Item(21, 4).Style
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '6.99'
$ws.Cells.Item(21, 4).Style = $st
$ws.Cells.Item(21, 5).Value = '  -1.01%  '
$st = $ws.Cells.Item(22, 4).Style
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '0.575'
$ws.Cells.Item(22, 4).Style = $st
$ws.Cells.Item(22, 5).Value = '  +7.34%  '
$st = $ws.Cells.Item(23, 4).Style
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '0.993'
$ws.Cells.Item(23, 4).Style = $st
$ws.Cells.Item(23, 5).Value = '  -0.68%  '
$st = $ws.Cells.Item(24, 4).Style
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '66.97'
$ws.Cells.Item(24, 4).Style = $st
$ws.Cells.Item(24, 5).Value = '  +1.54%  '
$st = $ws.Cells.Item(25, 4).Style
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '0.175'
$ws.Cells.Item(25, 4).Style = $st
$ws.Cells.Item(25, 5).Value = '  +4.59%  '
$st = $ws.Cells.Item(26, 4).Style
$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '8.75'
$ws.Cells.Item(26, 4).Style = $st
$ws.Cells.Item(26, 5).Value = '  +1.18%  '
$ws.Cells.Item(27, 4).Value = '0.0₃0953'
$ws.Cells.Item(27, 5).Value = '  +8.84%  '
$st = $ws.Cells.Item(28, 4).Style
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '0.997'
$ws.Cells.Item(28, 4).Style = $st
$ws.Cells.Item(28, 5).Value = '  +0.15%  '
$st = $ws.Cells.Item(29, 4).Style
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '2.04'
$ws.Cells.Item(29, 4).Style = $st
$ws.Cells.Item(29, 5).Value = '  -0.17%  '
$st = $ws.Cells.Item(30, 4).Style
$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '7.18'
$ws.Cells.Item(30, 4).Style = $st
$ws.Cells.Item(30, 5).Value = '  +0.23%  '
$st = $ws.Cells.Item(31, 4).Style
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '1.25'
$ws.Cells.Item(31, 4).Style = $st
$ws.Cells.Item(31, 5).Value = '  +4.98%  '
$st = $ws.Cells.Item(32, 4).Style
$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '170.91'
$ws.Cells.Item(32, 4).Style = $st
$ws.Cells.Item(32, 5).Value = '  +0.40%  '
$st = $ws.Cells.Item(33, 4).Style
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '5.09'
$ws.Cells.Item(33, 4).Style = $st
$ws.Cells.Item(33, 5).Value = '  +6.59%  '
$st = $ws.Cells.Item(34, 4).Style
$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '0.998'
$ws.Cells.Item(34, 4).Style = $st
$ws.Cells.Item(34, 5).Value = '  +0.17%  '
$st = $ws.Cells.Item(35, 4).Style
$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '20.74'
$ws.Cells.Item(35, 4).Style = $st
$ws.Cells.Item(35, 5).Value = '  +0.66%  '
$st = $ws.Cells.Item(36, 4).Style
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '1.48'
$ws.Cells.Item(36, 4).Style = $st
$ws.Cells.Item(36, 5).Value = '  +2.41%  '
$st = $ws.Cells.Item(37, 4).Style
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '1.84'
$ws.Cells.Item(37, 4).Style = $st
$ws.Cells.Item(37, 5).Value = '  +1.70%  '
$st = $ws.Cells.Item(38, 4).Style
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '1.02'
$ws.Cells.Item(38, 4).Style = $st
$ws.Cells.Item(38, 5).Value = '  -0.61%  '
$st = $ws.Cells.Item(39, 4).Style
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '4.24'
$ws.Cells.Item(39, 4).Style = $st
$ws.Cells.Item(39, 5).Value = '  -0.51%  '
$st = $ws.Cells.Item(40, 4).Style
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '337.24'
$ws.Cells.Item(40, 4).Style = $st
$ws.Cells.Item(40, 5).Value = '  -4.62%  '
$st = $ws.Cells.Item(41, 4).Style
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '6.21'
$ws.Cells.Item(41, 4).Style = $st
$ws.Cells.Item(41, 5).Value = '  +8.54%  '
$st = $ws.Cells.Item(42, 4).Style
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '40.08'
$ws.Cells.Item(42, 4).Style = $st
$ws.Cells.Item(42, 5).Value = '  +1.96%  '
$st = $ws.Cells.Item(43, 4).Style
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '22.31'
$ws.Cells.Item(43, 4).Style = $st
$ws.Cells.Item(43, 5).Value = '  +0.17%  '
$st = $ws.Cells.Item(44, 4).Style
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '0.0606'
$ws.Cells.Item(44, 4).Style = $st
$ws.Cells.Item(44, 5).Value = '  +2.04%  '
$st = $ws.Cells.Item(45, 4).Style
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '22.25'
$ws.Cells.Item(45, 4).Style = $st
$ws.Cells.Item(45, 5).Value = '  +0.45%  '
$st = $ws.Cells.Item(46, 4).Style
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '0.649'
$ws.Cells.Item(46, 4).Style = $st
$ws.Cells.Item(46, 5).Value = '  +1.04%  '
$st = $ws.Cells.Item(47, 4).Style
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '0.0261'
$ws.Cells.Item(47, 4).Style = $st
$ws.Cells.Item(47, 5).Value = '  +0.52%  '
$st = $ws.Cells.Item(48, 4).Style
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '137.66'
$ws.Cells.Item(48, 4).Style = $st
$ws.Cells.Item(48, 5).Value = '  -3.43%  '
$st = $ws.Cells.Item(49, 4).Style
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '0.103'
$ws.Cells.Item(49, 4).Style = $st
$ws.Cells.Item(49, 5).Value = '  +1.24%  '
$st = $ws.Cells.Item(50, 4).Style
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '0.999'
$ws.Cells.Item(50, 4).Style = $st
$ws.Cells.Item(50, 5).Value = '  +0.88%  '
$ws.Cells.Item(51, 4).Value = '2.151.18'
$ws.Cells.Item(51, 5).Value = '  -0.93%  '
